$p = $ppt.ActivePresentation

$p.Slides.Item(3).Shapes.Item(2).TextFrame.TextRange.Text = "5%"
$p.Slides.Item(5).Shapes.Item(2).TextFrame.TextRange.Text = "15%"
$p.Slides.Item(7).Shapes.Item(2).TextFrame.TextRange.Text = "50%"
$p.Slides.Item(8).Shapes.Item(2).TextFrame.TextRange.Text = "30%"
